$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 19
$ws.Range("I2").Value = 67
$ws.Range("J2").Value = 274
$ws.Range("L2").Value = 64
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 51
$ws.Range("P2").Value = 0
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 39
$ws.Range("T2").Value = 48
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 433
$ws.Range("X2").Value = 372
$ws.Range("Z2").Value = 4
$ws.Range("AA2").Value = 1
